$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the header row (row 1) and delete it, shifting the remaining
# data rows up. This removes the "DIV 1 / DIV 2 / DIV 3" header row.
$ws.Rows.Item(1).Select()
$ws.Rows.Item(1).Delete()
